$d = $word.ActiveDocument

# --- 1. Finish the incomplete "Problem Statement" paragraph ---

# 1a. "camera." -> "cameras." in "...field of view of the camera. An EMG..."
$d.Content.Find.Execute(
    "the field of view of the camera. An EMG",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "the field of view of the cameras. An EMG",
    2) | Out-Null

# 1b. Replace the "(THIS PARAGRAPH IS INCOMPLETE)" placeholder (and the two
#     spaces preceding it) with the real continuation of the paragraph.
$d.Content.Find.Execute(
    "(up to 300 feet).  (THIS PARAGRAPH IS INCOMPLETE)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "(up to 300 feet). Also, the processing required for gesture detection is significantly less than the image processing that is currently in use.   ",
    2) | Out-Null

# --- 2. Remove the stray leading comma before "Often" in the next paragraph ---
$d.Content.Find.Execute(
    "theft. , Often,",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "theft. Often,",
    2) | Out-Null

# --- 3. Bump the cached header DATE field result by one day ---
$hdr = $d.Sections(1).Headers(1)
$hdr.Range.Find.Execute(
    "October 3, 2017",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "October 4, 2017",
    2) | Out-Null

# --- 4. Add 6pt "space after" to the Subtitle paragraph style ---
$subtitle = $d.Styles("Subtitle")
$subtitle.ParagraphFormat.SpaceAfter = 6
